$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '328.59'
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '1.02%'
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '5.513'
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '0.37%'
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-0.39%'
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '1.982'
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '3.59%'
$ws.Range("E5").Style = "Normal"

$ws.Range("B6").Value = 'BTSEToken'

$ws.Range("C6").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '2.584'
$ws.Range("D6").Style = "Normal"

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-4.83%'
$ws.Range("E6").Style = "Normal"

$ws.Range("B7").Value = 'MXToken'

$ws.Range("C7").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9498'
$ws.Range("D7").Style = "Normal"

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '0.82%'
$ws.Range("E7").Style = "Normal"

$ws.Range("B8").Value = 'LiechtensteinCryptoassetsExchange'

$ws.Range("C8").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.1137'
$ws.Range("D8").Style = "Normal"

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-2.13%'
$ws.Range("E8").Style = "Normal"

$ws.Range("B9").Value = 'WazirX'

$ws.Range("C9").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1886'
$ws.Range("D9").Style = "Normal"

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '0.83%'
$ws.Range("E9").Style = "Normal"

$ws.Range("B10").Value = 'MCDex'

$ws.Range("C10").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '10.66'
$ws.Range("D10").Style = "Normal"

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '25.04%'
$ws.Range("E10").Style = "Normal"

$ws.Range("B11").Value = 'MandalaExchangeToken'

$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.09961'
$ws.Range("D11").Style = "Normal"

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-0.27%'
$ws.Range("E11").Style = "Normal"

$ws.Range("B12").Value = 'BitrueCoin'

$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.04713'
$ws.Range("D12").Style = "Normal"

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '11.41%'
$ws.Range("E12").Style = "Normal"

$ws.Range("B13").Value = 'BitMartToken'

$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.1065'
$ws.Range("D13").Style = "Normal"

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-0.05%'
$ws.Range("E13").Style = "Normal"

$ws.Range("B14").Value = 'BitForexToken'

$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001267'
$ws.Range("D14").Style = "Normal"

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-1.03%'
$ws.Range("E14").Style = "Normal"

$ws.Range("B15").Value = 'CoinExToken'

$ws.Range("C15").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.04072'
$ws.Range("D15").Style = "Normal"

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-4.04%'
$ws.Range("E15").Style = "Normal"

$ws.Range("B16").Value = 'TigerCash'

$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.005973'
$ws.Range("D16").Style = "Normal"

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '1.39%'
$ws.Range("E16").Style = "Normal"

$ws.Range("B17").Value = 'OKB'

$ws.Range("C17").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '43.97'
$ws.Range("D17").Style = "Normal"

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-1.25%'
$ws.Range("E17").Style = "Normal"

$ws.Range("B18").Value = 'LEO'

$ws.Range("C18").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.364'
$ws.Range("D18").Style = "Normal"

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-6.25%'
$ws.Range("E18").Style = "Normal"

$ws.Range("B19").Value = 'GateToken'

$ws.Range("C19").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.369'
$ws.Range("D19").Style = "Normal"

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '1.97%'
$ws.Range("E19").Style = "Normal"

$ws.Range("B20").Value = 'BitpandaEcosystemToken'

$ws.Range("C20").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.3473'
$ws.Range("D20").Style = "Normal"

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '-0.33%'
$ws.Range("E20").Style = "Normal"

$ws.Range("B21").Value = 'ProBitToken'

$ws.Range("C21").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1419'
$ws.Range("D21").Style = "Normal"

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '3.41%'
$ws.Range("E21").Style = "Normal"

$ws.Range("B22").Value = 'ZBToken'

$ws.Range("C22").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.2545'
$ws.Range("D22").Style = "Normal"

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '0.61%'
$ws.Range("E22").Style = "Normal"

$ws.Range("B23").Value = 'BitKan'

$ws.Range("C23").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.001264'
$ws.Range("D23").Style = "Normal"

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '2.20%'
$ws.Range("E23").Style = "Normal"

$ws.Range("B24").Value = 'HotbitToken'

$ws.Range("C24").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.004333'
$ws.Range("D24").Style = "Normal"

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-4.86%'
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0001200'
$ws.Range("D25").Style = "Normal"

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-0.02%'
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0003742'
$ws.Range("D26").Style = "Normal"

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02594'
$ws.Range("D38").Style = "Normal"

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '-1.69%'
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05678'
$ws.Range("D39").Style = "Normal"

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '3.74%'
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.007553'
$ws.Range("D40").Style = "Normal"

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-1.87%'
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1399'
$ws.Range("D41").Style = "Normal"

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '0.44%'
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.007484'
$ws.Range("D42").Style = "Normal"

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '4.25%'
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002016'
$ws.Range("D43").Style = "Normal"

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-1.57%'
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.008628'
$ws.Range("D44").Style = "Normal"

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-0.40%'
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00007100'
$ws.Range("D45").Style = "Normal"

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-0.11%'
$ws.Range("E45").Style = "Normal"

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-0.02%'
$ws.Range("E46").Style = "Normal"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '55.39%'
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.003770'
$ws.Range("D48").Style = "Normal"

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '4.54%'
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00002101'
$ws.Range("D49").Style = "Normal"

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-0.02%'
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0002001'
$ws.Range("D50").Style = "Normal"

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.02%'
$ws.Range("E50").Style = "Normal"
